# Refined metadata to be additional tab
#
# 1. Re-stamp the "panel_query_time" (F column) values on the "data" sheet
#    with the new query run's timestamps.
# 2. Add a new "metadata" worksheet (placed after "data") describing the
#    PanelApp query that produced the "data" sheet.

$wb = $excel.ActiveWorkbook
$dataWs = $wb.Worksheets.Item(1)
$dataWs.Name = "data"

# --- 1. Update the F column (time_taken) timestamps on the "data" sheet ---
$dataWs.Range("F2").Value = "2021-10-05 14:20:06.172628"
$dataWs.Range("F3").Value = "2021-10-05 14:20:06.172636"
$dataWs.Range("F4").Value = "2021-10-05 14:20:06.172639"
$dataWs.Range("F5").Value = "2021-10-05 14:20:06.172642"
$dataWs.Range("F6").Value = "2021-10-05 14:20:06.172645"
$dataWs.Range("F7").Value = "2021-10-05 14:20:06.172648"
$dataWs.Range("F8").Value = "2021-10-05 14:20:06.172651"
$dataWs.Range("F9").Value = "2021-10-05 14:20:06.172654"
$dataWs.Range("F10").Value = "2021-10-05 14:20:06.172656"
$dataWs.Range("F11").Value = "2021-10-05 14:20:06.172659"
$dataWs.Range("F12").Value = "2021-10-05 14:20:06.172662"
$dataWs.Range("F13").Value = "2021-10-05 14:20:06.172665"
$dataWs.Range("F14").Value = "2021-10-05 14:20:06.172667"
$dataWs.Range("F15").Value = "2021-10-05 14:20:06.172670"
$dataWs.Range("F16").Value = "2021-10-05 14:20:06.172673"
$dataWs.Range("F17").Value = "2021-10-05 14:20:06.172675"
$dataWs.Range("F18").Value = "2021-10-05 14:20:06.172678"
$dataWs.Range("F19").Value = "2021-10-05 14:20:06.172681"
$dataWs.Range("F20").Value = "2021-10-05 14:20:06.172683"
$dataWs.Range("F21").Value = "2021-10-05 14:20:06.172686"
$dataWs.Range("F22").Value = "2021-10-05 14:20:06.172689"
$dataWs.Range("F23").Value = "2021-10-05 14:20:06.172691"
$dataWs.Range("F24").Value = "2021-10-05 14:20:06.172694"
$dataWs.Range("F25").Value = "2021-10-05 14:20:06.172697"
$dataWs.Range("F26").Value = "2021-10-05 14:20:06.172700"
$dataWs.Range("F27").Value = "2021-10-05 14:20:06.172702"
$dataWs.Range("F28").Value = "2021-10-05 14:20:06.172705"
$dataWs.Range("F29").Value = "2021-10-05 14:20:06.172708"
$dataWs.Range("F30").Value = "2021-10-05 14:20:06.172710"
$dataWs.Range("F31").Value = "2021-10-05 14:20:06.172713"

# --- 2. Add the "metadata" sheet, right after "data" ---
$metaWs = $wb.Worksheets.Add($null, $dataWs)
$metaWs.Name = "metadata"

$metaWs.Range("B1").Value = "data_name"
$metaWs.Range("C1").Value = "data_id"
$metaWs.Range("D1").Value = "data_version"
$metaWs.Range("E1").Value = "data_version_created"
$metaWs.Range("F1").Value = "panel_query_time"
$metaWs.Range("G1").Value = "panel_get_request"

$metaWs.Range("A2").Value = 0
$metaWs.Range("B2").Value = "Ectodermal dysplasia without a known gene mutation"
$metaWs.Range("C2").Value = 136
# "data_version" is a text value ("1.22"), not a number - force a text
# number format before assigning so it is not coerced to a float.
$metaWs.Range("D2").NumberFormat = "@"
$metaWs.Range("D2").Value = "1.22"
$metaWs.Range("E2").Value = "2021-07-22T15:12:47.428047Z"
$metaWs.Range("F2").Value = "2021-10-05 14:20:06.169287"
$metaWs.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/136/?format=json"

# Match the header style used by the "data" sheet's header row (bold,
# bordered, center/top aligned) for row 1 of "metadata", and the same
# style for A2 (as used for the "data" sheet's leading index column).
$headerRange = $metaWs.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$idxRange = $metaWs.Range("A2")
$idxRange.Font.Bold = $true
$idxRange.Borders.LineStyle = 1
$idxRange.HorizontalAlignment = -4108
$idxRange.VerticalAlignment = -4160
